{"js": "// Replace each arithmetic-problem answer in the table with its updated\n// value. Every table cell holds exactly one run of text shaped like\n// \"39-6=33\"; we look each original formula up with Body.search() (an\n// exact, case-sensitive, non-wildcard match) and swap it for the new\n// formula with Range.insertText(..., Word.InsertLocation.replace). This\n// only rewrites the text inside the existing run, so the run/paragraph\n// formatting (font, size, alignment, etc.) is left untouched.\nconst replacements = [\n  [\"39-6=33\", \"13-6=7\"],\n  [\"85-63=22\", \"79-58=21\"],\n  [\"3+92=95\", \"10+87=97\"],\n  [\"89-78=11\", \"91-44=47\"],\n  [\"0+92=92\", \"34+43=77\"],\n  [\"91-62=29\", \"92-20=72\"],\n  [\"12+16=28\", \"31+34=65\"],\n  [\"38-26=12\", \"81-39=42\"],\n  [\"54-0=54\", \"90-54=36\"],\n  [\"37+20=57\", \"58-26=32\"],\n  [\"19+65=84\", \"15+73=88\"],\n  [\"47+39=86\", \"20+33=53\"],\n  [\"1+51=52\", \"51+1=52\"],\n  [\"94-71=23\", \"73-38=35\"],\n  [\"29+53=82\", \"14+2=16\"],\n  [\"23-0=23\", \"88-9=79\"],\n  [\"74-36=38\", \"65+18=83\"],\n  [\"35-2=33\", \"69-15=54\"],\n  [\"73-48=25\", \"46-3=43\"],\n  [\"98-37=61\", \"38+31=69\"],\n  [\"67-61=6\", \"6+60=66\"],\n  [\"67-62=5\", \"77-61=16\"],\n  [\"89+8=97\", \"81-33=48\"],\n  [\"1+24=25\", \"35+49=84\"],\n  [\"67-55=12\", \"10+33=43\"],\n  [\"35-12=23\", \"91-0=91\"],\n  [\"47+45=92\", \"67-15=52\"],\n  [\"23+4=27\", \"38-28=10\"],\n  [\"55-51=4\", \"18-18=0\"],\n  [\"66-35=31\", \"86-33=53\"],\n  [\"31+35=66\", \"29+70=99\"],\n  [\"13+13=26\", \"68+7=75\"],\n  [\"57+18=75\", \"32-12=20\"],\n  [\"95-6=89\", \"49+46=95\"],\n  [\"18+32=50\", \"71-9=62\"],\n  [\"57-14=43\", \"19+4=23\"],\n  [\"77+4=81\", \"90-13=77\"],\n  [\"20+50=70\", \"56+43=99\"],\n  [\"42+25=67\", \"95-90=5\"],\n  [\"17+59=76\", \"78+13=91\"],\n  [\"0+64=64\", \"18+66=84\"],\n  [\"46+10=56\", \"16+56=72\"],\n  [\"95-18=77\", \"26-21=5\"],\n  [\"87-29=58\", \"23+34=57\"],\n  [\"59-32=27\", \"69+18=87\"],\n  [\"76+3=79\", \"73-11=62\"],\n  [\"35-34=1\", \"33+61=94\"],\n  [\"49-47=2\", \"13+22=35\"],\n  [\"88-20=68\", \"88-3=85\"],\n  [\"36-31=5\", \"26+53=79\"],\n  [\"29-23=6\", \"14+75=89\"],\n  [\"35-27=8\", \"69-6=63\"],\n  [\"6+42=48\", \"1+74=75\"],\n  [\"60-31=29\", \"15+73=88\"],\n  [\"16+37=53\", \"89-17=72\"],\n  [\"75-20=55\", \"95-3=92\"],\n  [\"3+84=87\", \"45+2=47\"],\n  [\"64-33=31\", \"99-5=94\"],\n  [\"7+54=61\", \"39+52=91\"],\n  [\"37+42=79\", \"14+15=29\"],\n  [\"73+13=86\", \"30+40=70\"],\n  [\"94-63=31\", \"63+33=96\"],\n  [\"33+32=65\", \"56+14=70\"],\n  [\"27-4=23\", \"66-46=20\"],\n  [\"34-9=25\", \"59-0=59\"],\n  [\"54-7=47\", \"83-55=28\"],\n  [\"49+11=60\", \"69-15=54\"],\n  [\"89-36=53\", \"97-45=52\"],\n  [\"21-19=2\", \"71+0=71\"],\n  [\"20+41=61\", \"49-7=42\"],\n  [\"9+51=60\", \"24+70=94\"],\n  [\"79-39=40\", \"87+3=90\"],\n  [\"2+65=67\", \"69+5=74\"],\n  [\"11+45=56\", \"20+11=31\"],\n  [\"20-0=20\", \"0+4=4\"],\n  [\"30-14=16\", \"33-14=19\"],\n  [\"97-72=25\", \"40-12=28\"],\n  [\"85+2=87\", \"40-38=2\"],\n  [\"77-74=3\", \"40+2=42\"],\n  [\"92-78=14\", \"41+7=48\"],\n  [\"29+50=79\", \"8+7=15\"],\n  [\"24+51=75\", \"87-28=59\"],\n  [\"72-3=69\", \"39+2=41\"],\n  [\"60-56=4\", \"17+73=90\"],\n  [\"84-23=61\", \"17-8=9\"],\n  [\"62-35=27\", \"22+32=54\"],\n  [\"76-27=49\", \"90-48=42\"],\n  [\"28-4=24\", \"11+28=39\"],\n  [\"61-27=34\", \"21+55=76\"],\n  [\"13+19=32\", \"20+31=51\"],\n  [\"96+1=97\", \"16+19=35\"],\n  [\"42+52=94\", \"19+66=85\"],\n  [\"32+34=66\", \"13+14=27\"],\n  [\"44-2=42\", \"40+2=42\"],\n  [\"83-48=35\", \"8+57=65\"],\n  [\"22+63=85\", \"70-65=5\"],\n  [\"64+34=98\", \"92-29=63\"],\n  [\"52-10=42\", \"3+78=81\"],\n  [\"15+65=80\", \"91-22=69\"],\n  [\"52-11=41\", \"76-50=26\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update every arithmetic-problem answer cell in the table to its new\n# value. Each table cell holds exactly one run of text like \"39-6=33\";\n# Find/Replace (whole-text match, not a wildcard pattern) swaps it for the\n# new formula while leaving the surrounding run/paragraph formatting intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('39-6=33', '13-6=7'),\n    @('85-63=22', '79-58=21'),\n    @('3+92=95', '10+87=97'),\n    @('89-78=11', '91-44=47'),\n    @('0+92=92', '34+43=77'),\n    @('91-62=29', '92-20=72'),\n    @('12+16=28', '31+34=65'),\n    @('38-26=12', '81-39=42'),\n    @('54-0=54', '90-54=36'),\n    @('37+20=57', '58-26=32'),\n    @('19+65=84', '15+73=88'),\n    @('47+39=86', '20+33=53'),\n    @('1+51=52', '51+1=52'),\n    @('94-71=23', '73-38=35'),\n    @('29+53=82', '14+2=16'),\n    @('23-0=23', '88-9=79'),\n    @('74-36=38', '65+18=83'),\n    @('35-2=33', '69-15=54'),\n    @('73-48=25', '46-3=43'),\n    @('98-37=61', '38+31=69'),\n    @('67-61=6', '6+60=66'),\n    @('67-62=5', '77-61=16'),\n    @('89+8=97', '81-33=48'),\n    @('1+24=25', '35+49=84'),\n    @('67-55=12', '10+33=43'),\n    @('35-12=23', '91-0=91'),\n    @('47+45=92', '67-15=52'),\n    @('23+4=27', '38-28=10'),\n    @('55-51=4', '18-18=0'),\n    @('66-35=31', '86-33=53'),\n    @('31+35=66', '29+70=99'),\n    @('13+13=26', '68+7=75'),\n    @('57+18=75', '32-12=20'),\n    @('95-6=89', '49+46=95'),\n    @('18+32=50', '71-9=62'),\n    @('57-14=43', '19+4=23'),\n    @('77+4=81', '90-13=77'),\n    @('20+50=70', '56+43=99'),\n    @('42+25=67', '95-90=5'),\n    @('17+59=76', '78+13=91'),\n    @('0+64=64', '18+66=84'),\n    @('46+10=56', '16+56=72'),\n    @('95-18=77', '26-21=5'),\n    @('87-29=58', '23+34=57'),\n    @('59-32=27', '69+18=87'),\n    @('76+3=79', '73-11=62'),\n    @('35-34=1', '33+61=94'),\n    @('49-47=2', '13+22=35'),\n    @('88-20=68', '88-3=85'),\n    @('36-31=5', '26+53=79'),\n    @('29-23=6', '14+75=89'),\n    @('35-27=8', '69-6=63'),\n    @('6+42=48', '1+74=75'),\n    @('60-31=29', '15+73=88'),\n    @('16+37=53', '89-17=72'),\n    @('75-20=55', '95-3=92'),\n    @('3+84=87', '45+2=47'),\n    @('64-33=31', '99-5=94'),\n    @('7+54=61', '39+52=91'),\n    @('37+42=79', '14+15=29'),\n    @('73+13=86', '30+40=70'),\n    @('94-63=31', '63+33=96'),\n    @('33+32=65', '56+14=70'),\n    @('27-4=23', '66-46=20'),\n    @('34-9=25', '59-0=59'),\n    @('54-7=47', '83-55=28'),\n    @('49+11=60', '69-15=54'),\n    @('89-36=53', '97-45=52'),\n    @('21-19=2', '71+0=71'),\n    @('20+41=61', '49-7=42'),\n    @('9+51=60', '24+70=94'),\n    @('79-39=40', '87+3=90'),\n    @('2+65=67', '69+5=74'),\n    @('11+45=56', '20+11=31'),\n    @('20-0=20', '0+4=4'),\n    @('30-14=16', '33-14=19'),\n    @('97-72=25', '40-12=28'),\n    @('85+2=87', '40-38=2'),\n    @('77-74=3', '40+2=42'),\n    @('92-78=14', '41+7=48'),\n    @('29+50=79', '8+7=15'),\n    @('24+51=75', '87-28=59'),\n    @('72-3=69', '39+2=41'),\n    @('60-56=4', '17+73=90'),\n    @('84-23=61', '17-8=9'),\n    @('62-35=27', '22+32=54'),\n    @('76-27=49', '90-48=42'),\n    @('28-4=24', '11+28=39'),\n    @('61-27=34', '21+55=76'),\n    @('13+19=32', '20+31=51'),\n    @('96+1=97', '16+19=35'),\n    @('42+52=94', '19+66=85'),\n    @('32+34=66', '13+14=27'),\n    @('44-2=42', '40+2=42'),\n    @('83-48=35', '8+57=65'),\n    @('22+63=85', '70-65=5'),\n    @('64+34=98', '92-29=63'),\n    @('52-10=42', '3+78=81'),\n    @('15+65=80', '91-22=69'),\n    @('52-11=41', '76-50=26'),\n)\n\nforeach ($pair in $replacements) {\n    $find = $pair[0]\n    $replace = $pair[1]\n    $range = $d.Content\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # MatchCase=True and MatchWildcards=False keep this an exact literal match\n    # (the '+'/'-'/'=' in the formulas are plain text, not wildcard tokens);\n    # Replace=wdReplaceOne(2) only touches the single matching cell.\n    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null\n}\n"}
